$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.799.71'
$ws.Range('E2').Value = '  +0.20%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.810.38'
$ws.Range('E3').Value = '  +0.64%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.01'
$ws.Range('E4').Value = '  +0.72%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '601.94'
$ws.Range('E5').Value = '  +1.05%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.21'
$ws.Range('E6').Value = '  -0.48%  '

$ws.Range('E7').Value = '  +0.10%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.518'
$ws.Range('E8').Value = '  -0.12%  '

$ws.Range('E9').Value = '  +0.30%  '

$ws.Range('E10').Value = '  +0.98%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.35'
$ws.Range('E11').Value = '  +0.85%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000250'
$ws.Range('E12').Value = '  -1.01%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '36.03'
$ws.Range('E13').Value = '  +0.24%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.450.83'
$ws.Range('E14').Value = '  +0.71%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.822.90'
$ws.Range('E15').Value = '  +0.76%  '

$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '67.847.33'
$ws.Range('E16').Value = '  +0.34%  '

$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.47'
$ws.Range('E17').Value = '  -0.77%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.09'
$ws.Range('E18').Value = '  +0.72%  '

$ws.Range('E19').Value = '  +1.68%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '463.91'
$ws.Range('E20').Value = '  +0.94%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.82'
$ws.Range('E21').Value = '  -1.93%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.703'
$ws.Range('E22').Value = '  +0.98%  '

$ws.Range('E23').Value = '  -3.17%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '83.22'
$ws.Range('E24').Value = '  -0.15%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.09'
$ws.Range('E25').Value = '  +0.84%  '

$ws.Range('E26').Value = '  +0.19%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.01'
$ws.Range('E27').Value = '  +0.03%  '

$ws.Range('E28').Value = '  -0.07%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.962.21'
$ws.Range('E29').Value = '  +0.78%  '

$ws.Range('E30').Value = '  +0.17%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.42'
$ws.Range('E31').Value = '  +3.19%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.21'
$ws.Range('E32').Value = '  -0.22%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '29.42'
$ws.Range('E33').Value = '  -0.73%  '

$ws.Range('E34').Value = '  +0.12%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.07'
$ws.Range('E35').Value = '  -0.21%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0998'
$ws.Range('E36').Value = '  -0.19%  '

$ws.Range('E37').Value = '  +0.15%  '

$ws.Range('B38').Value = 'Mantle'
$ws.Range('C38').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.996'
$ws.Range('E38').Value = '  -0.08%  '

$ws.Range('B39').Value = 'Filecoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.81'
$ws.Range('E39').Value = '  +0.92%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.23'
$ws.Range('E40').Value = '  -4.02%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  +0.14%  '

$ws.Range('E42').Value = '  +0.02%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '45.16'
$ws.Range('E43').Value = '  -0.82%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '47.77'
$ws.Range('E44').Value = '  -0.76%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.299'
$ws.Range('E45').Value = '  +0.17%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '28.06'
$ws.Range('E46').Value = '  +5.16%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '151.24'
$ws.Range('E47').Value = '  +0.79%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.36'
$ws.Range('E48').Value = '  +0.54%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.37'
$ws.Range('E49').Value = '  +10.37%  '

$ws.Range('E50').Value = '  +1.78%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '391.45'
$ws.Range('E51').Value = '  -0.40%  '

